$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9789442420005798
$ws.Range("B1").Value = 1.344437479972839
$ws.Range("C1").Value = 1.946287274360657
$ws.Range("D1").Value = 5.308160781860352
$ws.Range("E1").Value = 1.873375296592712
